$wb = $excel.ActiveWorkbook

# "rotation list" sheet: remove row 6 (YYNENb), shifting former row 7 (YYENEr) up into row 6.
$ws1 = $wb.Worksheets.Item("rotation list")
$ws1.Rows.Item(6).Delete()

# "rotation con1" sheet: remove rows 9 and 10 (the YYNENb pair), shifting former rows 11-12
# (the YYENEz pair) up into rows 9-10.
$ws2 = $wb.Worksheets.Item("rotation con1")
$ws2.Rows.Item(9).Delete()
$ws2.Rows.Item(9).Delete()
